# Scheduled-runner update: refresh Universalis market-price snapshots and the
# resulting profit figures (currentAveragePrice / LevePrice / LeveProfit) on
# each job sheet. Values below are pre-computed by the external price puller;
# this script only needs to land them on the right cells.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 14007.315
$ws.Range("I21").Value = 10708.177
$ws.Range("J21").Value = 42050
$ws.Range("K21").Value = 10708.177
$ws.Range("L21").Value = 42050
$ws.Range("M21").Value = -10240.177
$ws.Range("N21").Value = -42986
$ws.Range("H23").Value = 14007.315
$ws.Range("I23").Value = 10708.177
$ws.Range("J23").Value = 42050
$ws.Range("K23").Value = 10708.177
$ws.Range("L23").Value = 42050
$ws.Range("M23").Value = -10474.177
$ws.Range("N23").Value = -42518
$ws.Range("H29").Value = 548.25
$ws.Range("I29").Value = 548.25
$ws.Range("K29").Value = 1644.75
$ws.Range("M29").Value = -1363.75
$ws.Range("H32").Value = 5003.3335
$ws.Range("I32").Value = 7445
$ws.Range("J32").Value = 3782.5
$ws.Range("K32").Value = 7445
$ws.Range("L32").Value = 3782.5
$ws.Range("M32").Value = -7119
$ws.Range("N32").Value = -4434.5
$ws.Range("H40").Value = 1026.2727
$ws.Range("I40").Value = 982.1111
$ws.Range("J40").Value = 1225
$ws.Range("K40").Value = 982.1111
$ws.Range("L40").Value = 1225
$ws.Range("M40").Value = -807.1111
$ws.Range("N40").Value = -1575
$ws.Range("H51").Value = 2983
$ws.Range("I51").Value = 2979
$ws.Range("J51").Value = 2984
$ws.Range("K51").Value = 2979
$ws.Range("L51").Value = 2984
$ws.Range("M51").Value = -2495
$ws.Range("N51").Value = -3952
$ws.Range("H62").Value = 2725.5
$ws.Range("I62").Value = 1766.6666
$ws.Range("J62").Value = 3045.111
$ws.Range("K62").Value = 1766.6666
$ws.Range("L62").Value = 3045.111
$ws.Range("M62").Value = -1142.6666
$ws.Range("N62").Value = -4293.111
$ws.Range("H65").Value = 2725.5
$ws.Range("I65").Value = 1766.6666
$ws.Range("J65").Value = 3045.111
$ws.Range("K65").Value = 8833.333000000001
$ws.Range("L65").Value = 15225.555
$ws.Range("M65").Value = -5713.333000000001
$ws.Range("N65").Value = -21465.555
$ws.Range("H113").Value = 1644
$ws.Range("I113").Value = 1610
$ws.Range("J113").Value = 1666.6666
$ws.Range("K113").Value = 1610
$ws.Range("L113").Value = 1666.6666
$ws.Range("M113").Value = 1644
$ws.Range("N113").Value = -8174.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1461.4445
$ws.Range("I102").Value = 1489.9375
$ws.Range("K102").Value = 1489.9375
$ws.Range("M102").Value = 132.0625
$ws.Range("H110").Value = 1020.1923
$ws.Range("I110").Value = 908.7273
$ws.Range("K110").Value = 908.7273
$ws.Range("M110").Value = 1136.2727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()   # column recalculated to blank

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2626.6667
$ws.Range("J49").Value = 2626.6667
$ws.Range("L49").Value = 7880.000100000001
$ws.Range("N49").Value = -8192.000100000001
$ws.Range("H105").Value = 181602820
$ws.Range("J105").Value = 181602820
$ws.Range("L105").Value = 544808460
$ws.Range("N105").Value = -544813702
$ws.Range("H114").Value = 1847.5
$ws.Range("J114").Value = 2779.1
$ws.Range("L114").Value = 8337.299999999999
$ws.Range("N114").Value = -14845.3
$ws.Range("H117").Value = 3599.8
$ws.Range("J117").Value = 4666.3335
$ws.Range("L117").Value = 13999.0005
$ws.Range("N117").Value = -20883.0005
$ws.Range("H129").Value = 729.1539
$ws.Range("I129").Value = 386.55554
$ws.Range("K129").Value = 1159.66662
$ws.Range("M129").Value = 3840.33338
$ws.Range("H131").Value = 2503556.5
$ws.Range("I131").Value = 5798.5454
$ws.Range("J131").Value = 3450982
$ws.Range("K131").Value = 17395.6362
$ws.Range("L131").Value = 10352946
$ws.Range("M131").Value = -12355.6362
$ws.Range("N131").Value = -10363026

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5286.6665
$ws.Range("I113").Value = 10980
$ws.Range("J113").Value = 2440
$ws.Range("K113").Value = 10980
$ws.Range("L113").Value = 2440
$ws.Range("M113").Value = -8810
$ws.Range("N113").Value = -6780
$ws.Range("H122").Value = 7152
$ws.Range("I122").Value = 8412
$ws.Range("J122").Value = 4002
$ws.Range("K122").Value = 25236
$ws.Range("L122").Value = 12006
$ws.Range("M122").Value = -22786
$ws.Range("N122").Value = -16906
$ws.Range("H132").Value = 3258.8333
$ws.Range("I132").Value = 2695.5833
$ws.Range("K132").Value = 8086.749899999999
$ws.Range("M132").Value = -5556.749899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2083.1924
$ws.Range("I7").Value = 1831.3334
$ws.Range("K7").Value = 1831.3334
$ws.Range("M7").Value = -1719.3334
$ws.Range("H40").Value = 2122.9333
$ws.Range("I40").Value = 2062.0833
$ws.Range("J40").Value = 2366.3333
$ws.Range("K40").Value = 2062.0833
$ws.Range("L40").Value = 2366.3333
$ws.Range("M40").Value = -1926.0833
$ws.Range("N40").Value = -2638.3333
$ws.Range("H82").Value = 1791.125
$ws.Range("I82").Value = 1638.25
$ws.Range("J82").Value = 2249.75
$ws.Range("K82").Value = 1638.25
$ws.Range("L82").Value = 2249.75
$ws.Range("M82").Value = -1277.25
$ws.Range("N82").Value = -2971.75
$ws.Range("H85").Value = 1791.125
$ws.Range("I85").Value = 1638.25
$ws.Range("J85").Value = 2249.75
$ws.Range("K85").Value = 1638.25
$ws.Range("L85").Value = 2249.75
$ws.Range("M85").Value = -390.25
$ws.Range("N85").Value = -4745.75
$ws.Range("H126").Value = 2083.1924
$ws.Range("I126").Value = 1831.3334
$ws.Range("K126").Value = 5494.0002
$ws.Range("M126").Value = -3024.0002
$ws.Range("H132").Value = 3116.2083
$ws.Range("I132").Value = 4026.0908
$ws.Range("J132").Value = 2346.3076
$ws.Range("K132").Value = 12078.2724
$ws.Range("L132").Value = 7038.9228
$ws.Range("M132").Value = -9548.2724
$ws.Range("N132").Value = -12098.9228
$ws.Range("H136").Value = 1771.7273
$ws.Range("I136").Value = 1158
$ws.Range("J136").Value = 2283.1667
$ws.Range("K136").Value = 3474
$ws.Range("L136").Value = 6849.500100000001
$ws.Range("M136").Value = -924
$ws.Range("N136").Value = -11949.5001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15981.667
$ws.Range("J54").Value = 15981.667
$ws.Range("L54").Value = 15981.667
$ws.Range("N54").Value = -17021.667
$ws.Range("H113").Value = 1069.1428
$ws.Range("I113").Value = 1033.6666
$ws.Range("K113").Value = 3100.9998
$ws.Range("M113").Value = -930.9998000000001
$ws.Range("H132").Value = 5620.467
$ws.Range("I132").Value = 9834
$ws.Range("J132").Value = 2811.4443
$ws.Range("K132").Value = 29502
$ws.Range("L132").Value = 8434.332900000001
$ws.Range("M132").Value = -26972
$ws.Range("N132").Value = -13494.3329
